# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF ---
# Copy the formatting of the existing last header cell (AC1) onto the
# new header cells so they match the bold/bordered/centered header style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-47): season record for every player is the same, ---
# --- since it reflects the team's overall record for the season.    ---
$wins = 97
$losses = 66
$ties = 0

for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
